$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 130
$ws.Range("A130").Value = "'2025-07-18"
$ws.Range("B130").Value = "Changchun Yatai"
$ws.Range("C130").Value = "SHANGHAI SIPG"
$ws.Range("D130").Value = 1
$ws.Range("E130").Value = 3
$ws.Range("F130").Value = 1341028
$ws.Range("G130").Value = 3
$ws.Range("H130").Value = 3
$ws.Range("I130").Value = 1
$ws.Range("J130").Value = 2
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = 0
$ws.Range("N130").Value = 0
$ws.Range("O130").Value = 1
$ws.Range("P130").Value = 3
$ws.Range("Q130").Value = "'39%"
$ws.Range("R130").Value = "'61%"
$ws.Range("S130").Value = "V"

# Row 131
$ws.Range("A131").Value = "'2025-07-18"
$ws.Range("B131").Value = "Wuhan Three Towns"
$ws.Range("C131").Value = "Qingdao Youth Island"
$ws.Range("D131").Value = 1
$ws.Range("E131").Value = 1
$ws.Range("F131").Value = 1341027
$ws.Range("G131").Value = 11
$ws.Range("H131").Value = 3
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 0
$ws.Range("N131").Value = 0
$ws.Range("O131").Value = 1
$ws.Range("P131").Value = 1
$ws.Range("Q131").Value = "'50%"
$ws.Range("R131").Value = "'50%"
$ws.Range("S131").Value = "E"

# Row 132
$ws.Range("A132").Value = "'2025-07-18"
$ws.Range("B132").Value = "Tianjin Teda"
$ws.Range("C132").Value = "Chengdu Better City"
$ws.Range("D132").Value = 2
$ws.Range("E132").Value = 1
$ws.Range("F132").Value = 1341029
$ws.Range("G132").Value = 3
$ws.Range("H132").Value = 3
$ws.Range("I132").Value = 1
$ws.Range("J132").Value = 5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 1
$ws.Range("M132").Value = 0
$ws.Range("N132").Value = 0
$ws.Range("O132").Value = 2
$ws.Range("P132").Value = 1
$ws.Range("Q132").Value = "'37%"
$ws.Range("R132").Value = "'63%"
$ws.Range("S132").Value = "L"

# Row 133
$ws.Range("A133").Value = "'2025-07-18"
$ws.Range("B133").Value = "Hangzhou Greentown"
$ws.Range("C133").Value = "Yunnan Yukun"
$ws.Range("D133").Value = 3
$ws.Range("E133").Value = 1
$ws.Range("F133").Value = 1341030
$ws.Range("G133").Value = 7
$ws.Range("H133").Value = 8
$ws.Range("I133").Value = 4
$ws.Range("J133").Value = 1
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = 0
$ws.Range("N133").Value = 0
$ws.Range("O133").Value = 3
$ws.Range("P133").Value = 1
$ws.Range("Q133").Value = "'54%"
$ws.Range("R133").Value = "'46%"
$ws.Range("S133").Value = "L"

# Row 134
$ws.Range("A134").Value = "'2025-07-19"
$ws.Range("B134").Value = "Dalian Zhixing"
$ws.Range("C134").Value = "Shandong Luneng"
$ws.Range("D134").Value = 2
$ws.Range("E134").Value = 0
$ws.Range("F134").Value = 1341031
$ws.Range("G134").Value = 7
$ws.Range("H134").Value = 8
$ws.Range("I134").Value = 2
$ws.Range("J134").Value = 1
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 0
$ws.Range("N134").Value = 0
$ws.Range("O134").Value = 2
$ws.Range("P134").Value = 0
$ws.Range("Q134").Value = "'42%"
$ws.Range("R134").Value = "'58%"
$ws.Range("S134").Value = "L"

# Row 135
$ws.Range("A135").Value = "'2025-07-19"
$ws.Range("B135").Value = "Beijing Guoan"
$ws.Range("C135").Value = "Shanghai Shenhua"
$ws.Range("D135").Value = 1
$ws.Range("E135").Value = 3
$ws.Range("F135").Value = 1341033
$ws.Range("G135").Value = 8
$ws.Range("H135").Value = 2
$ws.Range("I135").Value = 1
$ws.Range("J135").Value = 3
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = 0
$ws.Range("N135").Value = 0
$ws.Range("O135").Value = 1
$ws.Range("P135").Value = 3
$ws.Range("Q135").Value = "'72%"
$ws.Range("R135").Value = "'28%"
$ws.Range("S135").Value = "V"

# Row 136
$ws.Range("A136").Value = "'2025-07-19"
$ws.Range("B136").Value = "Sichuan Jiuniu"
$ws.Range("C136").Value = "Qingdao Jonoon"
$ws.Range("D136").Value = 4
$ws.Range("E136").Value = 0
$ws.Range("F136").Value = 1341032
$ws.Range("G136").Value = 5
$ws.Range("H136").Value = 7
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 1
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 0
$ws.Range("N136").Value = 0
$ws.Range("O136").Value = 4
$ws.Range("P136").Value = 0
$ws.Range("Q136").Value = "'46%"
$ws.Range("R136").Value = "'54%"
$ws.Range("S136").Value = "L"

# Row 137
$ws.Range("A137").Value = "'2025-07-19"
$ws.Range("B137").Value = "Henan Jianye"
$ws.Range("C137").Value = "Meizhou Kejia"
$ws.Range("D137").Value = 1
$ws.Range("E137").Value = 1
$ws.Range("F137").Value = 1341034
$ws.Range("G137").Value = 9
$ws.Range("H137").Value = 5
$ws.Range("I137").Value = 2
$ws.Range("J137").Value = 1
$ws.Range("K137").Value = 1
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 0
$ws.Range("N137").Value = 0
$ws.Range("O137").Value = 1
$ws.Range("P137").Value = 1
$ws.Range("Q137").Value = "'51%"
$ws.Range("R137").Value = "'49%"
$ws.Range("S137").Value = "E"
